# The published SectorGroup.xlsx codelist re-orders its trailing columns so
# that the human-readable "group"/"category" names lead and the numeric
# codes trail:
#
#   before: code | name | status | category-code | group-code | group-name | category-name
#   after:  code | name | status | group-name    | category-name | group-code | category-code
#
# i.e. for every row (including the header row) column D swaps with column F,
# and column E swaps with column G.
#
# We do this with Range.Copy (not .Value assignment) so the original cell
# type is preserved exactly: the numeric-looking codes ("111", "110", ...)
# stay shared-string text cells instead of being reinterpreted as numbers,
# and no new cell styles get introduced along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$srcRange = "D1:G" + $lastRow

# Stage the current D:G block into scratch columns I:L.
$ws.Range($srcRange).Copy($ws.Range("I1"))
# I=old D (category-code), J=old E (group-code), K=old F (group-name), L=old G (category-name)

# Write the swapped columns back from the scratch copy.
$ws.Range("K1:K" + $lastRow).Copy($ws.Range("D1"))   # D = group-name
$ws.Range("L1:L" + $lastRow).Copy($ws.Range("E1"))   # E = category-name
$ws.Range("J1:J" + $lastRow).Copy($ws.Range("F1"))   # F = group-code
$ws.Range("I1:I" + $lastRow).Copy($ws.Range("G1"))   # G = category-code

# Remove the scratch columns.
$ws.Range("I1:L" + $lastRow).Clear()
